$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append 10 new rows (147-156) of Mac-Address test data, following the
# same pattern as the existing rows: regcntr_id, device_id, lang_code,
# is_active, cr_by, cr_dtimes
$startRow = 147
$startDevice = 3000166

for ($i = 0; $i -lt 10; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = 10001
    $ws.Cells.Item($r, 2).Value = $startDevice + $i
    $ws.Cells.Item($r, 3).Value = "eng"
    $ws.Cells.Item($r, 4).Value = $true
    $ws.Cells.Item($r, 5).Value = "superadmin"
    $ws.Cells.Item($r, 6).Value = "now()"
}

# Update the visible selection to match the author's final cursor
# position after entering the new data.
$null = $ws.Range("C152").Select()
